$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.691.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.316.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.75%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "269.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.62%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.64%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0934"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.61%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.91%  "

# Row 13
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.674.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.14%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.88%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.853"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.328.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.79%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.689.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000106"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.92%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.28%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.19%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.85%  "

# Row 29
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.33%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0889"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.67%  "

# Row 36
$ws.Range("E36").Value = "  +0.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0348"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.98%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.235"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +15.71%  "

# Row 41
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.28%  "

# Row 42
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.53%  "

# Row 43
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +19.23%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
$ws.Range("B45").Value = "MultiversX"
$ws.Range("C45").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.52%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.78%  "

# Row 47
$ws.Range("E47").Value = "  +2.82%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49
$ws.Range("E49").Value = "  -0.16%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.555.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.429"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.49%  "
